# "updated code for 2015": refresh the latest-year count and append the
# new 2015 row to the impact-country-numbers-per-year table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2014's count was revised from 428 to 432.
$ws.Range("B5").Value = 432

# Append the new 2015 row. A6 must land as a text shared-string (like the
# other project_year labels), not a number, so force text formatting
# before assigning the value, then drop the number format back to the
# sheet's default style so no extra formatting sticks to the cell.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2015"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = 292
